$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "299.51"
Set-TextValue $ws.Range("E2") "0.13%"
$ws.Range("F2").Value = "18-1-2023"
Set-TextValue $ws.Range("G2") "0"

Set-TextValue $ws.Range("D3") "32.13"
Set-TextValue $ws.Range("E3") "1.78%"
$ws.Range("F3").Value = "18-1-2023"
Set-TextValue $ws.Range("G3") "0"

Set-TextValue $ws.Range("D4") "4.975"
Set-TextValue $ws.Range("E4") "-3.09%"
$ws.Range("F4").Value = "18-1-2023"
Set-TextValue $ws.Range("G4") "0"

Set-TextValue $ws.Range("D5") "0.07868"
Set-TextValue $ws.Range("E5") "-1.51%"
$ws.Range("F5").Value = "18-1-2023"
Set-TextValue $ws.Range("G5") "0"

Set-TextValue $ws.Range("D6") "2.222"
Set-TextValue $ws.Range("E6") "-11.09%"
$ws.Range("F6").Value = "18-1-2023"
Set-TextValue $ws.Range("G6") "0"

Set-TextValue $ws.Range("D7") "7.793"
Set-TextValue $ws.Range("E7") "-0.10%"
$ws.Range("F7").Value = "18-1-2023"
Set-TextValue $ws.Range("G7") "0"

Set-TextValue $ws.Range("D8") "3.845"
Set-TextValue $ws.Range("E8") "-1.52%"
$ws.Range("F8").Value = "18-1-2023"
Set-TextValue $ws.Range("G8") "0"

Set-TextValue $ws.Range("D9") "0.9215"
Set-TextValue $ws.Range("E9") "-0.60%"
$ws.Range("F9").Value = "18-1-2023"
Set-TextValue $ws.Range("G9") "0"

Set-TextValue $ws.Range("E10") "-1.00%"
$ws.Range("F10").Value = "18-1-2023"
Set-TextValue $ws.Range("G10") "0"

Set-TextValue $ws.Range("D11") "0.07929"
Set-TextValue $ws.Range("E11") "6.90%"
$ws.Range("F11").Value = "18-1-2023"
Set-TextValue $ws.Range("G11") "0"

Set-TextValue $ws.Range("D12") "0.08593"
Set-TextValue $ws.Range("E12") "-3.52%"
$ws.Range("F12").Value = "18-1-2023"
Set-TextValue $ws.Range("G12") "0"

Set-TextValue $ws.Range("D13") "0.03083"
Set-TextValue $ws.Range("E13") "1.55%"
$ws.Range("F13").Value = "18-1-2023"
Set-TextValue $ws.Range("G13") "0"

Set-TextValue $ws.Range("E14") "-0.12%"
$ws.Range("F14").Value = "18-1-2023"
Set-TextValue $ws.Range("G14") "0"

Set-TextValue $ws.Range("D15") "0.001523"
Set-TextValue $ws.Range("E15") "0.74%"
$ws.Range("F15").Value = "18-1-2023"
Set-TextValue $ws.Range("G15") "0"

Set-TextValue $ws.Range("D16") "0.005954"
Set-TextValue $ws.Range("E16") "-3.31%"
$ws.Range("F16").Value = "18-1-2023"
Set-TextValue $ws.Range("G16") "0"

Set-TextValue $ws.Range("E17") "2,100.52%"
$ws.Range("F17").Value = "18-1-2023"
Set-TextValue $ws.Range("G17") "0"

Set-TextValue $ws.Range("E18") "-1.16%"
$ws.Range("F18").Value = "18-1-2023"
Set-TextValue $ws.Range("G18") "0"

Set-TextValue $ws.Range("D19") "2.177"
Set-TextValue $ws.Range("E19") "-4.97%"
$ws.Range("F19").Value = "18-1-2023"
Set-TextValue $ws.Range("G19") "0"

$ws.Range("F20").Value = "18-1-2023"
Set-TextValue $ws.Range("G20") "0"

Set-TextValue $ws.Range("D21") "0.1282"
Set-TextValue $ws.Range("E21") "-2.77%"
$ws.Range("F21").Value = "18-1-2023"
Set-TextValue $ws.Range("G21") "0"

Set-TextValue $ws.Range("D22") "4.298"
Set-TextValue $ws.Range("E22") "2.79%"
$ws.Range("F22").Value = "18-1-2023"
Set-TextValue $ws.Range("G22") "0"

Set-TextValue $ws.Range("D23") "0.1794"
Set-TextValue $ws.Range("E23") "5.05%"
$ws.Range("F23").Value = "18-1-2023"
Set-TextValue $ws.Range("G23") "0"

Set-TextValue $ws.Range("D24") "0.04594"
Set-TextValue $ws.Range("E24") "-0.52%"
$ws.Range("F24").Value = "18-1-2023"
Set-TextValue $ws.Range("G24") "0"

Set-TextValue $ws.Range("D25") "0.001229"
Set-TextValue $ws.Range("E25") "-1.01%"
$ws.Range("F25").Value = "18-1-2023"
Set-TextValue $ws.Range("G25") "0"

Set-TextValue $ws.Range("D26") "0.004413"
Set-TextValue $ws.Range("E26") "-2.09%"
$ws.Range("F26").Value = "18-1-2023"
Set-TextValue $ws.Range("G26") "0"

Set-TextValue $ws.Range("D27") "0.0001250"
Set-TextValue $ws.Range("E27") "3.96%"
$ws.Range("F27").Value = "18-1-2023"
Set-TextValue $ws.Range("G27") "0"

$ws.Range("F28").Value = "18-1-2023"
Set-TextValue $ws.Range("G28") "0"

$ws.Range("F29").Value = "18-1-2023"
Set-TextValue $ws.Range("G29") "0"

$ws.Range("F30").Value = "18-1-2023"
Set-TextValue $ws.Range("G30") "0"

$ws.Range("F31").Value = "18-1-2023"
Set-TextValue $ws.Range("G31") "0"

$ws.Range("F32").Value = "18-1-2023"
Set-TextValue $ws.Range("G32") "0"

$ws.Range("F33").Value = "18-1-2023"
Set-TextValue $ws.Range("G33") "0"

$ws.Range("F34").Value = "18-1-2023"
Set-TextValue $ws.Range("G34") "0"

$ws.Range("F35").Value = "18-1-2023"
Set-TextValue $ws.Range("G35") "0"

$ws.Range("F36").Value = "18-1-2023"
Set-TextValue $ws.Range("G36") "0"

$ws.Range("F37").Value = "18-1-2023"
Set-TextValue $ws.Range("G37") "0"

$ws.Range("F38").Value = "18-1-2023"
Set-TextValue $ws.Range("G38") "0"

Set-TextValue $ws.Range("D39") "0.01738"
Set-TextValue $ws.Range("E39") "-1.30%"
$ws.Range("F39").Value = "18-1-2023"
Set-TextValue $ws.Range("G39") "0"

Set-TextValue $ws.Range("D40") "0.04759"
Set-TextValue $ws.Range("E40") "3.39%"
$ws.Range("F40").Value = "18-1-2023"
Set-TextValue $ws.Range("G40") "0"

Set-TextValue $ws.Range("D41") "0.007503"
Set-TextValue $ws.Range("E41") "7.95%"
$ws.Range("F41").Value = "18-1-2023"
Set-TextValue $ws.Range("G41") "0"

Set-TextValue $ws.Range("D42") "0.1350"
Set-TextValue $ws.Range("E42") "-1.50%"
$ws.Range("F42").Value = "18-1-2023"
Set-TextValue $ws.Range("G42") "0"

Set-TextValue $ws.Range("D43") "0.002361"
Set-TextValue $ws.Range("E43") "10.87%"
$ws.Range("F43").Value = "18-1-2023"
Set-TextValue $ws.Range("G43") "0"

Set-TextValue $ws.Range("D44") "0.01182"
Set-TextValue $ws.Range("E44") "14.35%"
$ws.Range("F44").Value = "18-1-2023"
Set-TextValue $ws.Range("G44") "0"

Set-TextValue $ws.Range("D45") "0.00005974"
Set-TextValue $ws.Range("E45") "-5.55%"
$ws.Range("F45").Value = "18-1-2023"
Set-TextValue $ws.Range("G45") "0"

Set-TextValue $ws.Range("D46") "0.00000000750"
Set-TextValue $ws.Range("E46") "0.06%"
$ws.Range("F46").Value = "18-1-2023"
Set-TextValue $ws.Range("G46") "0"

Set-TextValue $ws.Range("D47") "0.003389"
Set-TextValue $ws.Range("E47") "-57.62%"
$ws.Range("F47").Value = "18-1-2023"
Set-TextValue $ws.Range("G47") "0"

Set-TextValue $ws.Range("D48") "0.8204"
$ws.Range("F48").Value = "18-1-2023"
Set-TextValue $ws.Range("G48") "0"

Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "0.06%"
$ws.Range("F49").Value = "18-1-2023"
Set-TextValue $ws.Range("G49") "0"

Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "0.06%"
$ws.Range("F50").Value = "18-1-2023"
Set-TextValue $ws.Range("G50") "0"

$ws.Range("F51").Value = "18-1-2023"
Set-TextValue $ws.Range("G51") "0"

Write-Output "done"